# Locate the end of the paragraph that ends with "them." (end of the
# "07.06.2021" diary entry) and insert two new BodyText paragraphs right
# after it:
#   1. an empty paragraph
#   2. a paragraph containing "I already had NodeJS installed"
# The pre-existing empty paragraph that originally followed "them." is
# left untouched and simply ends up after the two new paragraphs.

$d = $word.ActiveDocument

$rng = $d.Content
$found = $rng.Find.Execute("them.", $true, $false, $false, $false, $false, `
                            $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not find anchor text 'them.' in the document."
}

# Collapse to the point right after "them." (still inside its paragraph,
# before its paragraph mark).
$rng.Collapse(0)

# Insert a paragraph break followed immediately by the new sentence; this
# creates the blank paragraph first, then the run of text, both of which
# inherit the "BodyText" paragraph style and "en-US" language from the
# surrounding text, exactly like the rest of the document.
$rng.InsertAfter([char]13 + "I already had NodeJS installed")

# Collapse to the end of the text we just inserted and add the paragraph
# break that separates our new text paragraph from the paragraph that
# follows it (the pre-existing empty paragraph).
$rng.Collapse(0)
$rng.InsertParagraphAfter()
